$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRange, $val) {
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $val
    $cellRange.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "22.015.51"
Set-TextValue $ws.Range("E2") "  -1.40%  "

Set-TextValue $ws.Range("D3") "1.550.94"
Set-TextValue $ws.Range("E3") "  -0.88%  "

Set-TextValue $ws.Range("D4") "1.000"
Set-TextValue $ws.Range("E4") "  -0.13%  "

Set-TextValue $ws.Range("E5") "  -0.06%  "

Set-TextValue $ws.Range("D6") "288.04"
Set-TextValue $ws.Range("E6") "  +0.40%  "

Set-TextValue $ws.Range("D7") "0.3922"
Set-TextValue $ws.Range("E7") "  +4.21%  "

Set-TextValue $ws.Range("D8") "0.3201"
Set-TextValue $ws.Range("E8") "  -1.98%  "

Set-TextValue $ws.Range("D9") "41.93"
Set-TextValue $ws.Range("E9") "  -8.03%  "

Set-TextValue $ws.Range("D10") "0.07260"
Set-TextValue $ws.Range("E10") "  -1.98%  "

Set-TextValue $ws.Range("E11") "  -4.11%  "

Set-TextValue $ws.Range("D12") "1.001"
Set-TextValue $ws.Range("E12") "  -0.15%  "

Set-TextValue $ws.Range("D13") "18.89"
Set-TextValue $ws.Range("E13") "  -7.38%  "

Set-TextValue $ws.Range("D14") "5.613"
Set-TextValue $ws.Range("E14") "  -3.99%  "

Set-TextValue $ws.Range("D15") "6.636"
Set-TextValue $ws.Range("E15") "  -2.45%  "

Set-TextValue $ws.Range("D16") "0.00001123"
Set-TextValue $ws.Range("E16") "  +2.93%  "

Set-TextValue $ws.Range("D17") "1.550.12"
Set-TextValue $ws.Range("E17") "  -1.33%  "

Set-TextValue $ws.Range("D18") "0.06572"
Set-TextValue $ws.Range("E18") "  -2.25%  "

Set-TextValue $ws.Range("D19") "83.55"
Set-TextValue $ws.Range("E19") "  -2.67%  "

Set-TextValue $ws.Range("E20") "  -0.12%  "

Set-TextValue $ws.Range("D21") "6.279"
Set-TextValue $ws.Range("E21") "  -1.44%  "

Set-TextValue $ws.Range("D22") "15.68"
Set-TextValue $ws.Range("E22") "  -3.52%  "

Set-TextValue $ws.Range("D23") "11.19"
Set-TextValue $ws.Range("E23") "  -4.22%  "

Set-TextValue $ws.Range("D24") "22.026.78"
Set-TextValue $ws.Range("E24") "  -1.37%  "

Set-TextValue $ws.Range("D25") "2.362"
Set-TextValue $ws.Range("E25") "  +2.98%  "

Set-TextValue $ws.Range("D26") "2.407"
Set-TextValue $ws.Range("E26") "  -4.32%  "

Set-TextValue $ws.Range("D27") "147.62"
Set-TextValue $ws.Range("E27") "  -1.90%  "

Set-TextValue $ws.Range("D28") "18.55"
Set-TextValue $ws.Range("E28") "  -4.67%  "

Set-TextValue $ws.Range("D29") "4.837"
Set-TextValue $ws.Range("E29") "  -1.36%  "

Set-TextValue $ws.Range("D30") "1.725.64"
Set-TextValue $ws.Range("E30") "  -1.23%  "

Set-TextValue $ws.Range("D31") "118.50"
Set-TextValue $ws.Range("E31") "  -4.05%  "

Set-TextValue $ws.Range("D32") "1.051"
Set-TextValue $ws.Range("E32") "  +0.88%  "

Set-TextValue $ws.Range("D33") "5.678"
Set-TextValue $ws.Range("E33") "  -3.94%  "

Set-TextValue $ws.Range("D34") "0.08341"
Set-TextValue $ws.Range("E34") "  +1.23%  "

Set-TextValue $ws.Range("D35") "9.105"
Set-TextValue $ws.Range("E35") "  -3.72%  "

Set-TextValue $ws.Range("E36") "  -16.46%  "

Set-TextValue $ws.Range("D37") "0.06147"
Set-TextValue $ws.Range("E37") "  -2.31%  "

Set-TextValue $ws.Range("D38") "0.02260"
Set-TextValue $ws.Range("E38") "  -4.90%  "

Set-TextValue $ws.Range("D39") "5.095"
Set-TextValue $ws.Range("E39") "  -3.13%  "

Set-TextValue $ws.Range("D40") "1.211"
Set-TextValue $ws.Range("E40") "  -5.03%  "

Set-TextValue $ws.Range("D41") "0.2062"
Set-TextValue $ws.Range("E41") "  -5.44%  "

Set-TextValue $ws.Range("E42") "  +0.06%  "

Set-TextValue $ws.Range("D43") "10.59"
Set-TextValue $ws.Range("E43") "  -3.96%  "

Set-TextValue $ws.Range("D44") "0.5793"
Set-TextValue $ws.Range("E44") "  -4.51%  "

Set-TextValue $ws.Range("D45") "13.23"
Set-TextValue $ws.Range("E45") "  -3.76%  "

Set-TextValue $ws.Range("D46") "3.732"
Set-TextValue $ws.Range("E46") "  -0.37%  "

Set-TextValue $ws.Range("D47") "0.5556"
Set-TextValue $ws.Range("E47") "  -5.53%  "

Set-TextValue $ws.Range("D48") "117.79"
Set-TextValue $ws.Range("E48") "  -4.67%  "

Set-TextValue $ws.Range("D49") "1.891"
Set-TextValue $ws.Range("E49") "  -5.29%  "

Set-TextValue $ws.Range("D50") "1.135"
Set-TextValue $ws.Range("E50") "  -3.60%  "

Set-TextValue $ws.Range("D51") "0.06819"
Set-TextValue $ws.Range("E51") "  -4.37%  "
